# Grilles_Eval.xlsx - "Update Bloc 3 - Intro"
# Adds a new worksheet "v2" (placed after "Feuil1") containing the intro /
# header block of a new evaluation grid layout, and adjusts the selection
# state left on "Feuil1".

$wb = $excel.ActiveWorkbook
$feuil1 = $wb.Worksheets.Item(1)

# --- Feuil1: leave the selection parked on the header row instead of the
#     bottom of the grid -------------------------------------------------
$feuil1.Range("D6:X6").Select()

# --- add the new sheet "v2" right after "Feuil1" ------------------------
$ws = $wb.Worksheets.Add($null, $feuil1)
$ws.Name = "v2"

# Column widths (character units) matching the grid of 3 "blocs" of 4
# narrow columns separated by thin spacer columns.
$ws.Columns.Item(1).ColumnWidth = 5
$ws.Columns.Item(2).ColumnWidth = 43.44140625
$ws.Columns.Item(3).ColumnWidth = 1.6640625
$ws.Range("D1:G1").EntireColumn.ColumnWidth = 2.77734375
$ws.Columns.Item(8).ColumnWidth = 1.6640625
$ws.Range("I1:L1").EntireColumn.ColumnWidth = 2.77734375
$ws.Columns.Item(13).ColumnWidth = 1.6640625
$ws.Range("N1:Q1").EntireColumn.ColumnWidth = 2.77734375
$ws.Columns.Item(18).ColumnWidth = 1.88671875

# --- Row 1 : title --------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 18
$ws.Range("A1").Value = "ONIP / Outils Numériques pour l'Ingénieur.e en Physique "
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 14

# --- Row 3 : "NOMS : " fill-in line ---------------------------------------
$ws.Range("A3").Value = "NOMS : "
$ws.Range("A3").Borders.Item(7).LineStyle = 1
$ws.Range("A3").Borders.Item(8).LineStyle = 1
$ws.Range("A3").Borders.Item(9).LineStyle = 1
$ws.Range("B3:Q3").Borders.Item(8).LineStyle = 1
$ws.Range("B3:Q3").Borders.Item(9).LineStyle = 1
$ws.Range("Q3").Borders.Item(10).LineStyle = 1

# --- Row 4 : "Gpe : " fill-in box, plus +/- headers for the 3 blocs ------
$ws.Range("A4").Value = "Gpe : "
$ws.Range("A4").Borders.Item(7).LineStyle = 1
$ws.Range("A4").Borders.Item(8).LineStyle = 1
$ws.Range("A4").Borders.Item(9).LineStyle = 1
$ws.Range("B4").Borders.Item(8).LineStyle = 1
$ws.Range("B4").Borders.Item(9).LineStyle = 1
$ws.Range("B4").Borders.Item(10).LineStyle = 1

$ws.Range("D4,G4,I4,L4,N4,Q4").HorizontalAlignment = -4108
$ws.Range("E4:F4,J4:K4,O4:P4").HorizontalAlignment = -4108
$ws.Range("D4").Value = "+"
$ws.Range("G4").Value = "-"
$ws.Range("I4").Value = "+"
$ws.Range("L4").Value = "-"
$ws.Range("N4").Value = "+"
$ws.Range("Q4").Value = "-"

# --- Row 5 : "BLOC 1 / BLOC 2 / BLOC 3" merged headers --------------------
$ws.Range("D5:G5,I5:L5,N5:Q5").HorizontalAlignment = -4108
$ws.Range("D5,I5,N5").Borders.Item(7).LineStyle = 1
$ws.Range("D5:Q5").Borders.Item(8).LineStyle = 1
$ws.Range("D5:Q5").Borders.Item(9).LineStyle = 1
$ws.Range("G5,L5,Q5").Borders.Item(10).LineStyle = 1

$ws.Range("D5").Value = "BLOC 1"
$ws.Range("I5").Value = "BLOC 2"
$ws.Range("N5").Value = "BLOC 3"

$ws.Range("D5:G5").Merge()
$ws.Range("I5:L5").Merge()
$ws.Range("N5:Q5").Merge()

# --- Row 6 : first criterion block title ----------------------------------
$ws.Range("A6").Value = "S'APPROPRIER"
$ws.Range("A6").Font.Bold = $true

# --- Row 7 : first evaluated criterion + its three answer boxes ----------
$ws.Range("B7").Value = "Compréhension du sujet"
$ws.Range("D7:G7,I7:L7,N7:Q7").Borders.LineStyle = 1

# --- Row 9 : second occurrence of the same criterion label ---------------
$ws.Range("B9").Value = "Compréhension du sujet"

# --- final view state: select row 8 and make "v2" the active sheet -------
$ws.Rows.Item(8).Select()
$ws.Activate()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
